# AKA changes to sdCard to be saved while merged with Johns
#
# Adds a third "AQ32 + LA Timing" results table (columns H:L, rows 10-17)
# to Sheet1, mirroring the existing "AQ32 + MTK NMEA" table in columns A:E.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row (row 10), columns H:L -------------------------------------
# Mirror the look of the existing header row 1 (bold + centered), and give
# the new table its own title in J10 (merged J10:K10, already merged in the
# sheet).
$ws.Range("H10").Value = "Hz"
$ws.Range("I10").Value = "max uSec"
$ws.Range("J10").Value = "AQ32 + LA Timing"

$ws.Range("H1:I1").Copy()
$ws.Range("H10:I10").PasteSpecial(-4122)
$ws.Range("J1:L1").Copy()
$ws.Range("J10:L10").PasteSpecial(-4122)

# Re-apply the text that PasteSpecial(formats) above did not touch.
$ws.Range("H10").Value = "Hz"
$ws.Range("I10").Value = "max uSec"
$ws.Range("J10").Value = "AQ32 + LA Timing"

# --- Data rows 11-17, columns H:L ------------------------------------------
$ws.Range("H11").Value = 1000
$ws.Range("H12").Value = 500
$ws.Range("H13").Value = 100
$ws.Range("H14").Value = 50
$ws.Range("H15").Value = 10
$ws.Range("H16").Value = 5
$ws.Range("H17").Value = 1

$ws.Range("I11:I17").Formula = "=1/H11*1000000"

$ws.Range("J12").Value = 66
$ws.Range("J13").Value = 13
$ws.Range("J14").Value = 4.5
$ws.Range("J15").Value = 66
$ws.Range("J16").Value = 1.5

$ws.Range("K11:K17").Formula = "=J11/I11"
$ws.Range("K11:K17").NumberFormat = "0.00%"

$ws.Range("L11").Formula = "=SUM(K11:K17)"
$ws.Range("L11").NumberFormat = "0.00%"

# Note seen by the CLI while collecting the row 15 (10 Hz) sample.
$ws.Range("L15").Value = "199 when doing CLI output"

# Leave the same cell selected as in the saved workbook.
$ws.Range("L15").Select() | Out-Null
